# Update leveling/crafting price data across all profession sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect refreshed market prices.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 58
$ws.Range("H58").Value = 309.15384

# Row 112
$ws.Range("H112").Value = 4020.7273
$ws.Range("J112").Value = 4416.25
$ws.Range("L112").Value = 13248.75
$ws.Range("N112").Value = -15464.75

# Row 137
$ws.Range("H137").Value = 2980.875
$ws.Range("I137").Value = 2474.5
$ws.Range("K137").Value = 7423.5
$ws.Range("M137").Value = -4873.5

# Row 138
$ws.Range("H138").Value = 2476.9648
$ws.Range("J138").Value = 3910.8845
$ws.Range("L138").Value = 11732.6535
$ws.Range("N138").Value = -22012.6535

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 117
$ws.Range("H117").Value = 42500
$ws.Range("J117").Value = 42500
$ws.Range("L117").Value = 42500
$ws.Range("N117").Value = -51678

# Row 122
$ws.Range("H122").Value = 1719.5
$ws.Range("I122").Value = 1585
$ws.Range("K122").Value = 4755
$ws.Range("M122").Value = -2305

# Row 132
$ws.Range("H132").Value = 2367.449
$ws.Range("I132").Value = 2354.2708
$ws.Range("K132").Value = 7062.812399999999
$ws.Range("M132").Value = -4532.812399999999

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# Row 2
$ws.Range("H2").Value = 62833
$ws.Range("J2").Value = 69249.5
$ws.Range("L2").Value = 69249.5
$ws.Range("N2").Value = -69475.5

# Row 94
$ws.Range("H94").Value = 1679.0312
$ws.Range("I94").Value = 1324.3
$ws.Range("K94").Value = 1324.3
$ws.Range("M94").Value = -873.3

# Row 134
$ws.Range("H134").Value = 4663.069
$ws.Range("I134").Value = 4651.0356
$ws.Range("K134").Value = 13953.1068
$ws.Range("M134").Value = -11418.1068

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 4993.857
$ws.Range("I16").Value = 3928.125
$ws.Range("K16").Value = 3928.125
$ws.Range("M16").Value = -3641.125

# Row 31
$ws.Range("H31").Value = 4944.2964
$ws.Range("I31").Value = 3505.4443
$ws.Range("K31").Value = 3505.4443
$ws.Range("M31").Value = -3210.4443

# Row 34
$ws.Range("H34").Value = 4944.2964
$ws.Range("I34").Value = 3505.4443
$ws.Range("K34").Value = 3505.4443
$ws.Range("M34").Value = -3303.4443

# Row 99
$ws.Range("H99").Value = 5299.8335
$ws.Range("I99").Value = 4859.9
$ws.Range("J99").Value = 7499.5
$ws.Range("K99").Value = 4859.9
$ws.Range("L99").Value = 7499.5
$ws.Range("M99").Value = -3361.9
$ws.Range("N99").Value = -10495.5

# Row 113
$ws.Range("H113").Value = 4993.857
$ws.Range("I113").Value = 3928.125
$ws.Range("K113").Value = 3928.125
$ws.Range("M113").Value = -1758.125

# Row 126
$ws.Range("H126").Value = 5299.8335
$ws.Range("I126").Value = 4859.9
$ws.Range("J126").Value = 7499.5
$ws.Range("K126").Value = 14579.7
$ws.Range("L126").Value = 22498.5
$ws.Range("M126").Value = -12109.7
$ws.Range("N126").Value = -27438.5

# Row 132
$ws.Range("H132").Value = 3421.1482
$ws.Range("I132").Value = 3798.739
$ws.Range("K132").Value = 11396.217
$ws.Range("M132").Value = -8866.217000000001

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 100.13333
$ws.Range("I2").Value = 18.25
$ws.Range("J2").Value = 193.71428
$ws.Range("K2").Value = 109.5
$ws.Range("L2").Value = 1162.28568
$ws.Range("M2").Value = 3.5
$ws.Range("N2").Value = -1388.28568

# Row 32
$ws.Range("H32").Value = 1260
$ws.Range("I32").Value = 1650
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 4950
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -4667
$ws.Range("N32").Value = -3566

# Row 46
$ws.Range("H46").Value = 449.25

# Row 107
$ws.Range("H107").Value = 732.6667
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").Value = ""   # clear (was 1440)

# Row 122
$ws.Range("H122").Value = 858.5714
$ws.Range("I122").Value = 497
$ws.Range("J122").Value = 918.8333
$ws.Range("K122").Value = 4473
$ws.Range("L122").Value = 8269.4997
$ws.Range("M122").Value = -2023
$ws.Range("N122").Value = -13169.4997

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# Row 58
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = ""   # clear (was -40554)

# Row 62
$ws.Range("H62").Value = 149999
$ws.Range("J62").Value = 149999
$ws.Range("L62").Value = 149999
$ws.Range("N62").Value = -151371

# Row 65
$ws.Range("H65").Value = 149999
$ws.Range("J65").Value = 149999
$ws.Range("L65").Value = 449997
$ws.Range("N65").Value = -456861

# Row 69
$ws.Range("H69").Value = 38599.6
$ws.Range("J69").Value = 38599.6
$ws.Range("L69").Value = 38599.6
$ws.Range("N69").Value = -40097.6

# Row 72
$ws.Range("H72").Value = 38599.6
$ws.Range("J72").Value = 38599.6
$ws.Range("L72").Value = 115798.8
$ws.Range("N72").Value = -123286.8

# Row 80
$ws.Range("H80").Value = 9422.333000000001
$ws.Range("J80").Value = 9625
$ws.Range("L80").Value = 9625
$ws.Range("N80").Value = -11621

# Row 83
$ws.Range("H83").Value = 9422.333000000001
$ws.Range("J83").Value = 9625
$ws.Range("L83").Value = 48125
$ws.Range("N83").Value = -58109

# Row 99
$ws.Range("H99").Value = 10273.143
$ws.Range("I99").Value = 9485.333000000001
$ws.Range("K99").Value = 9485.333000000001
$ws.Range("M99").Value = -7239.333000000001

# Row 113
$ws.Range("H113").Value = 238270.77
$ws.Range("I113").Value = 365898.53
$ws.Range("K113").Value = 365898.53
$ws.Range("M113").Value = -363728.53

# Row 126
$ws.Range("H126").Value = 5466.6665
$ws.Range("I126").Value = 5466.6665
$ws.Range("K126").Value = 16399.9995
$ws.Range("M126").Value = -13929.9995

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 4199
$ws.Range("I46").Value = 4199
$ws.Range("K46").Value = 4199
$ws.Range("M46").Value = -4011

# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").Value = ""   # clear (was -47251)

# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").Value = ""   # clear (was -140256)

# Row 82
$ws.Range("H82").Value = 1950
$ws.Range("J82").Value = 2300
$ws.Range("L82").Value = 2300
$ws.Range("N82").Value = -3022

# Row 85
$ws.Range("H85").Value = 1950
$ws.Range("J85").Value = 2300
$ws.Range("L85").Value = 2300
$ws.Range("N85").Value = -4796

# Row 93
$ws.Range("H93").Value = 8839.581
$ws.Range("I93").Value = 2427.0476
$ws.Range("J93").Value = 22305.9
$ws.Range("K93").Value = 2427.0476
$ws.Range("L93").Value = 22305.9
$ws.Range("M93").Value = -1179.0476
$ws.Range("N93").Value = -24801.9

# Row 121
$ws.Range("H121").Value = 49999
$ws.Range("J121").Value = 49999
$ws.Range("L121").Value = 49999
$ws.Range("N121").Value = -53493

# Row 132
$ws.Range("H132").Value = 9775.115
$ws.Range("I132").Value = 8742.200000000001
$ws.Range("K132").Value = 26226.6
$ws.Range("M132").Value = -23696.6

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 36073.668
$ws.Range("I4").Value = 36073.668
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 36073.668
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -35960.668
$ws.Range("N4").Value = ""   # clear (was -49226)

# Row 110
$ws.Range("H110").Value = 72742
$ws.Range("I110").Value = 70585
$ws.Range("J110").Value = 74899
$ws.Range("K110").Value = 70585
$ws.Range("L110").Value = 74899
$ws.Range("M110").Value = -66495
$ws.Range("N110").Value = -83079

# Row 121
$ws.Range("H121").Value = 528959.5
$ws.Range("J121").Value = 528959.5
$ws.Range("L121").Value = 528959.5
$ws.Range("N121").Value = -532453.5

# Row 132
$ws.Range("H132").Value = 2812.4893
$ws.Range("I132").Value = 3021.946
$ws.Range("J132").Value = 2037.5
$ws.Range("K132").Value = 9065.838
$ws.Range("L132").Value = 6112.5
$ws.Range("M132").Value = -6535.838
$ws.Range("N132").Value = -11172.5

# Row 136
$ws.Range("H136").Value = 4474.9062
$ws.Range("I136").Value = 2175.125
$ws.Range("K136").Value = 6525.375
$ws.Range("M136").Value = -3975.375

Write-Output "applied edits"